$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$m.Shapes.Item(1).TextFrame.TextRange.Text = "Click to edit Master title style"
